# "tugas 1,2,3 menu user"
# The kategori_id column (A) is removed: kategori_kode (previously column B)
# becomes the new column A, and kategori_nama (previously column C) becomes
# the new column B. The old column C is cleared entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift kategori_kode (col B) -> col A, kategori_nama (col C) -> col B.
# Value2 is used because bulk multi-cell reads via Value are unreliable here.
$ws.Range("A1:A6").Value2 = $ws.Range("B1:B6").Value2
$ws.Range("B1:B6").Value2 = $ws.Range("C1:C6").Value2

# Fully remove the now-redundant third column's cell contents/styles.
$ws.Range("C1:C6").Clear()

# Widen column B now that it holds the longer "kategori_kode"/"kategori_nama"
# header instead of the old best-fit "kategori_kode" code values.
$ws.Columns.Item(2).ColumnWidth = 16.08

# Reflect the new active selection left behind after the edit.
$ws.Range("C4").Select() | Out-Null
